$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.490.25'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.727.34'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('D4').Value = "'0.9993"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'245.72"
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('D6').Value = "'0.9996"
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = "'0.4799"
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('D8').Value = "'0.2686"
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '1.728.22'
$ws.Range('D11').Value = "'0.07138"
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('D13').Value = "'0.6198"
$ws.Range('E13').Value = '  +5.30%  '
$ws.Range('D14').Value = "'4.518"
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('D15').Value = "'77.20"
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = "'0.9995"
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '26.508.28'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = "'0.9997"
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = "'0.000006953"
$ws.Range('E19').Value = '  +2.12%  '
$ws.Range('D20').Value = "'11.68"
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = '1.950.73'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('D22').Value = "'4.536"
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = "'8.950"
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('D24').Value = "'5.296"
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').Value = "'136.40"
$ws.Range('D26').Value = "'15.36"
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('D27').Value = "'1.804"
$ws.Range('E27').Value = '  +2.33%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').Value = "'106.76"
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = "'3.975"
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = "'0.08030"
$ws.Range('E31').Value = '  +3.63%  '
$ws.Range('D32').Value = "'3.725"
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D33').Value = "'0.04567"
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('D34').Value = "'0.9992"
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('D36').Value = "'0.6369"
$ws.Range('E36').Value = '  +2.54%  '
$ws.Range('D37').Value = "'0.9911"
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('D38').Value = "'0.9392"
$ws.Range('E38').Value = '  +2.85%  '
$ws.Range('D39').Value = "'2.098"
$ws.Range('E39').Value = '  +10.21%  '
$ws.Range('D40').Value = "'2.412"
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('D41').Value = "'104.71"
$ws.Range('E41').Value = '  -6.70%  '
$ws.Range('D42').Value = "'1.006"
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').Value = "'5.706"
$ws.Range('E43').Value = '  +9.25%  '
$ws.Range('D44').Value = "'0.01502"
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').Value = "'0.3915"
$ws.Range('E45').Value = '  +2.80%  '
$ws.Range('D46').Value = "'6.963"
$ws.Range('E46').Value = '  +11.61%  '
$ws.Range('D47').Value = "'0.1190"
$ws.Range('E47').Value = '  +4.02%  '
$ws.Range('D48').Value = "'0.05326"
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').Value = "'31.06"
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('D50').Value = "'7.885"
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('D51').Value = "'1.268"
$ws.Range('E51').Value = '  +4.02%  '
